# Populate the "ZEUS" (col J) and "%Difference_Best" (col K) results that
# were computed but never written back into the comparison sheet.
# Row 4 already had (stale) values that get overwritten with the final run's
# numbers; every other data row (3, 5-38, 40, 42, 43) gets the pair added.
# Rows 39, 41 and 44 are intentionally left alone (no ZEUS run recorded).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(3, 10).Value = 722.67529296875
$ws.Cells.Item(3, 11).Value = 37.0
$ws.Cells.Item(4, 10).Value = 1762.634521484375
$ws.Cells.Item(4, 11).Value = 33.0
$ws.Cells.Item(5, 10).Value = 722.67529296875
$ws.Cells.Item(5, 11).Value = 37.0
$ws.Cells.Item(6, 10).Value = 1133.586669921875
$ws.Cells.Item(6, 11).Value = 35.0
$ws.Cells.Item(7, 10).Value = 2895.0810546875
$ws.Cells.Item(7, 11).Value = 42.0
$ws.Cells.Item(8, 10).Value = 1133.586669921875
$ws.Cells.Item(8, 11).Value = 35.0
$ws.Cells.Item(9, 10).Value = 1149.4837646484375
$ws.Cells.Item(9, 11).Value = 39.0
$ws.Cells.Item(10, 10).Value = 2680.65234375
$ws.Cells.Item(10, 11).Value = 31.0
$ws.Cells.Item(11, 10).Value = 1149.4837646484375
$ws.Cells.Item(11, 11).Value = 39.0
$ws.Cells.Item(12, 10).Value = 2286.556396484375
$ws.Cells.Item(12, 11).Value = 43.0
$ws.Cells.Item(13, 10).Value = 1258.539794921875
$ws.Cells.Item(13, 11).Value = 61.0
$ws.Cells.Item(14, 10).Value = 1764.2425537109375
$ws.Cells.Item(14, 11).Value = 47.0
$ws.Cells.Item(15, 10).Value = 6050.01123046875
$ws.Cells.Item(15, 11).Value = 72.0
$ws.Cells.Item(16, 10).Value = 954.8070678710938
$ws.Cells.Item(16, 11).Value = 0.0
$ws.Cells.Item(17, 10).Value = 1862.6314697265625
$ws.Cells.Item(17, 11).Value = 0.0
$ws.Cells.Item(18, 10).Value = 2875.24072265625
$ws.Cells.Item(18, 11).Value = 0.0
$ws.Cells.Item(19, 10).Value = 1703.8721923828125
$ws.Cells.Item(19, 11).Value = 6.0
$ws.Cells.Item(20, 10).Value = 3215.44482421875
$ws.Cells.Item(20, 11).Value = 2.0
$ws.Cells.Item(21, 10).Value = 4846.4912109375
$ws.Cells.Item(21, 11).Value = 0.0
$ws.Cells.Item(22, 10).Value = 8367.3974609375
$ws.Cells.Item(22, 11).Value = 0.0
$ws.Cells.Item(23, 10).Value = 2476.455322265625
$ws.Cells.Item(23, 11).Value = 14.0
$ws.Cells.Item(24, 10).Value = 4583.64013671875
$ws.Cells.Item(24, 11).Value = 9.0
$ws.Cells.Item(25, 10).Value = 6795.0166015625
$ws.Cells.Item(25, 11).Value = 5.0
$ws.Cells.Item(26, 10).Value = 5045.349609375
$ws.Cells.Item(26, 11).Value = 36.0
$ws.Cells.Item(27, 10).Value = 5045.349609375
$ws.Cells.Item(27, 11).Value = 33.0
$ws.Cells.Item(28, 10).Value = 5040.25439453125
$ws.Cells.Item(28, 11).Value = 32.0
$ws.Cells.Item(29, 10).Value = 33524.46484375
$ws.Cells.Item(29, 11).Value = 52.0
$ws.Cells.Item(30, 10).Value = 33297.10546875
$ws.Cells.Item(30, 11).Value = 49.0
$ws.Cells.Item(31, 10).Value = 33343.578125
$ws.Cells.Item(31, 11).Value = 47.0
$ws.Cells.Item(32, 10).Value = 124680.8046875
$ws.Cells.Item(32, 11).Value = 67.0
$ws.Cells.Item(33, 10).Value = 125849.21875
$ws.Cells.Item(33, 11).Value = 64.0
$ws.Cells.Item(34, 10).Value = 124782.9140625
$ws.Cells.Item(34, 11).Value = 59.0
$ws.Cells.Item(35, 10).Value = 2734.0546875
$ws.Cells.Item(35, 11).Value = 23.0
$ws.Cells.Item(36, 10).Value = 5192.5439453125
$ws.Cells.Item(36, 11).Value = 37.0
$ws.Cells.Item(37, 10).Value = 8074.99609375
$ws.Cells.Item(37, 11).Value = 56.0
$ws.Cells.Item(38, 10).Value = 9463.0966796875
$ws.Cells.Item(38, 11).Value = 59.0
$ws.Cells.Item(40, 10).Value = 13173.0625
$ws.Cells.Item(40, 11).Value = 59.0
$ws.Cells.Item(42, 10).Value = 10987.5703125
$ws.Cells.Item(42, 11).Value = 57.0
$ws.Cells.Item(43, 10).Value = 16791.888671875
$ws.Cells.Item(43, 11).Value = 66.0
